$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: Houston
$ws.Range("A2").Value = "Houston"
$ws.Range("B2").Value = "https://www.sports-reference.com/cbb/schools/houston/2023.html"

# Update row 3: UCONN
$ws.Range("A3").Value = "UCONN"
$ws.Range("B3").Value = "https://www.sports-reference.com/cbb/schools/connecticut/2023.html"

# Delete row 4 (previously Wake Forest) entirely
$ws.Rows.Item(4).Delete()

# Update selection to A4
$ws.Range("A4").Select()
